# Remove two empty paragraphs from the document:
#   1. The empty, centered (sz=20) paragraph that sits right after the
#      "(Фамилия, имя, отчество)" caption and before the "Дата рождения"
#      table.
#   2. The empty, justified (sz=18) paragraph that sits right after the
#      closing diagnosis/recommendations table and before the
#      "Председатель подкомиссии" signature block.
#
# Deletions are applied starting with the one that occurs later in the
# document first, so that the character offsets used for the earlier
# deletion are not invalidated by the later one.

$d = $word.ActiveDocument

# --- 2) Empty paragraph right after the last table, before the
#        signatures block ("Председатель подкомиссии ...") ---
$tableCount = $d.Tables.Count
$lastTable = $d.Tables.Item($tableCount)
$afterTableStart = $lastTable.Range.End
$emptyParaAfterTable = $d.Range($afterTableStart, $afterTableStart + 1)
$emptyParaAfterTable.Delete()

# --- 1) Empty paragraph right after "(Фамилия, имя, отчество)",
#        before the "Дата рождения" table ---
$findRange = $d.Content
$found = $findRange.Find.Execute("(Фамилия, имя, отчество)", $false, $false,
                                  $false, $false, $false, $true, 1, $false,
                                  "", 0)
if ($found) {
    $findRange.Collapse(0)
    $emptyParaAfterCaption = $findRange.Next(4, 1)
    $emptyParaAfterCaption.Delete()
}
